# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# for a batch of Leve rows across all eight crafting-job sheets, per the
# scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 331.2
$ws.Range("I28").Value = 299.6111
$ws.Range("J28").Value = 615.5
$ws.Range("K28").Value = 299.6111
$ws.Range("L28").Value = 615.5
$ws.Range("M28").Value = 185.3889
$ws.Range("N28").Value = -1585.5
# Row 108 (Leve Item ID 25638)
$ws.Range("H108").Value = 31215.6
$ws.Range("J108").Value = 31215.6
$ws.Range("L108").Value = 31215.6
$ws.Range("N108").Value = -38895.6
# Row 109 (Leve Item ID 25639)
$ws.Range("H109").Value = 26735.2
$ws.Range("J109").Value = 26735.2
$ws.Range("L109").Value = 26735.2
$ws.Range("N109").Value = -29509.2
# Row 114 (Leve Item ID 25959)
$ws.Range("H114").Value = 39699.332
$ws.Range("J114").Value = 39699.332
$ws.Range("L114").Value = 39699.332
$ws.Range("N114").Value = -48377.332
# Row 117 (Leve Item ID 26118)
$ws.Range("H117").Value = 34064
$ws.Range("J117").Value = 34064
$ws.Range("L117").Value = 34064
$ws.Range("N117").Value = -43242
# Row 126 (Leve Item ID 34391)
$ws.Range("H126").Value = 46951.2
$ws.Range("J126").Value = 46951.2
$ws.Range("L126").Value = 46951.2
$ws.Range("N126").Value = -56831.2
# Row 128 (Leve Item ID 34540)
$ws.Range("H128").Value = 44169.168
$ws.Range("J128").Value = 44169.168
$ws.Range("L128").Value = 44169.168
$ws.Range("N128").Value = -54129.168
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 2485.0435
$ws.Range("J129").Value = 2991.923
$ws.Range("L129").Value = 8975.769
$ws.Range("N129").Value = -18975.769
# Row 130 (Leve Item ID 34691)
$ws.Range("H130").Value = 45401.6
$ws.Range("J130").Value = 45401.6
$ws.Range("L130").Value = 45401.6
$ws.Range("N130").Value = -55441.6
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 3825.1277
$ws.Range("I137").Value = 1022.4167
$ws.Range("K137").Value = 3067.2501
$ws.Range("M137").Value = -517.2501000000002

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 80 (Leve Item ID 10667)
$ws.Range("H80").Value = 51625.668
$ws.Range("J80").Value = 51625.668
$ws.Range("L80").Value = 51625.668
$ws.Range("N80").Value = -53621.668
# Row 83 (Leve Item ID 10667)
$ws.Range("H83").Value = 51625.668
$ws.Range("J83").Value = 51625.668
$ws.Range("L83").Value = 154877.004
$ws.Range("N83").Value = -164861.004
# Row 109 (Leve Item ID 25646)
$ws.Range("H109").Value = 40120.332
$ws.Range("J109").Value = 40120.332
$ws.Range("L109").Value = 40120.332
$ws.Range("N109").Value = -42894.332
# Row 118 (Leve Item ID 26150)
$ws.Range("H118").Value = 31270
$ws.Range("J118").Value = 31270
$ws.Range("L118").Value = 31270
$ws.Range("N118").Value = -34584
# Row 123 (Leve Item ID 34107)
$ws.Range("H123").Value = 35610.5
$ws.Range("J123").Value = 35610.5
$ws.Range("L123").Value = 35610.5
$ws.Range("N123").Value = -45410.5
# Row 124 (Leve Item ID 34252)
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""  # cell cleared (removed from XML)
# Row 125 (Leve Item ID 34251)
$ws.Range("H125").Value = 49469
$ws.Range("J125").Value = 49469
$ws.Range("L125").Value = 49469
$ws.Range("N125").Value = -59309
# Row 128 (Leve Item ID 34570)
$ws.Range("H128").Value = 49875
$ws.Range("J128").Value = 49875
$ws.Range("L128").Value = 49875
$ws.Range("N128").Value = -59835
# Row 130 (Leve Item ID 34732)
$ws.Range("H130").Value = 43664
$ws.Range("J130").Value = 43664
$ws.Range("L130").Value = 43664
$ws.Range("N130").Value = -53704
# Row 131 (Leve Item ID 34706)
$ws.Range("H131").Value = 51609
$ws.Range("J131").Value = 51609
$ws.Range("L131").Value = 51609
$ws.Range("N131").Value = -61689

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 117 (Leve Item ID 26124)
$ws.Range("H117").Value = 44999
$ws.Range("J117").Value = 44999
$ws.Range("L117").Value = 44999
$ws.Range("N117").Value = -54177
# Row 124 (Leve Item ID 34245)
$ws.Range("H124").Value = 47992
$ws.Range("J124").Value = 47992
$ws.Range("L124").Value = 47992
$ws.Range("N124").Value = -57812
# Row 125 (Leve Item ID 34235)
$ws.Range("H125").Value = 50772
$ws.Range("J125").Value = 50772
$ws.Range("L125").Value = 50772
$ws.Range("N125").Value = -60612
# Row 126 (Leve Item ID 34398)
$ws.Range("H126").Value = 44005.332
$ws.Range("J126").Value = 44005.332
$ws.Range("L126").Value = 44005.332
$ws.Range("N126").Value = -53885.332
# Row 130 (Leve Item ID 34682)
$ws.Range("H130").Value = 48273.75
$ws.Range("J130").Value = 48273.75
$ws.Range("L130").Value = 48273.75
$ws.Range("N130").Value = -58313.75

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 20 (Leve Item ID 34533)
$ws.Range("H20").Value = 49321.2
$ws.Range("J20").Value = 49321.2
$ws.Range("L20").Value = 49321.2
$ws.Range("N20").Value = -49793.2
# Row 30 (Leve Item ID 34533)
$ws.Range("H30").Value = 49321.2
$ws.Range("J30").Value = 49321.2
$ws.Range("L30").Value = 49321.2
$ws.Range("N30").Value = -49503.2
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 167164.28
$ws.Range("I31").Value = 1442.4193
$ws.Range("K31").Value = 1442.4193
$ws.Range("M31").Value = -1147.4193
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 167164.28
$ws.Range("I34").Value = 1442.4193
$ws.Range("K34").Value = 1442.4193
$ws.Range("M34").Value = -1240.4193
# Row 100 (Leve Item ID 34388)
$ws.Range("H100").Value = 33385
$ws.Range("J100").Value = 47770
$ws.Range("L100").Value = 47770
$ws.Range("N100").Value = -49934
# Row 128 (Leve Item ID 34533)
$ws.Range("H128").Value = 49321.2
$ws.Range("J128").Value = 49321.2
$ws.Range("L128").Value = 49321.2
$ws.Range("N128").Value = -59281.2

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 4501.2905
$ws.Range("J131").Value = 1841.9474
$ws.Range("L131").Value = 5525.8422
$ws.Range("N131").Value = -15605.8422

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 45941.285
$ws.Range("J130").Value = 45941.285
$ws.Range("L130").Value = 45941.285
$ws.Range("N130").Value = -55981.285

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 36 (Leve Item ID 34261)
$ws.Range("H36").Value = 46807.332
$ws.Range("J36").Value = 46807.332
$ws.Range("L36").Value = 46807.332
$ws.Range("N36").Value = -47931.332
# Row 119 (Leve Item ID 26288)
$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084
# Row 120 (Leve Item ID 26311)
$ws.Range("H120").Value = 51188.4
$ws.Range("J120").Value = 51188.4
$ws.Range("L120").Value = 51188.4
$ws.Range("N120").Value = -60864.4
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2355.8462
$ws.Range("I122").Value = 2355.8462
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7067.5386
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4617.5386
$ws.Range("N122").Value = ""  # cell cleared (removed from XML)
# Row 127 (Leve Item ID 34401)
$ws.Range("H127").Value = 50707.5
$ws.Range("J127").Value = 50707.5
$ws.Range("L127").Value = 50707.5
$ws.Range("N127").Value = -60627.5
# Row 130 (Leve Item ID 34729)
$ws.Range("H130").Value = 37996
$ws.Range("J130").Value = 37996
$ws.Range("L130").Value = 37996
$ws.Range("N130").Value = -48036

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 16 (Leve Item ID 26304)
$ws.Range("H16").Value = 45944.25
$ws.Range("J16").Value = 45944.25
$ws.Range("L16").Value = 45944.25
$ws.Range("N16").Value = -46528.25
# Row 110 (Leve Item ID 25825)
$ws.Range("H110").Value = 25433.6
$ws.Range("J110").Value = 25433.6
$ws.Range("L110").Value = 25433.6
$ws.Range("N110").Value = -33613.6
# Row 119 (Leve Item ID 26289)
$ws.Range("H119").Value = 36845
$ws.Range("J119").Value = 36845
$ws.Range("L119").Value = 36845
$ws.Range("N119").Value = -46521
# Row 120 (Leve Item ID 26310)
$ws.Range("H120").Value = 35206
$ws.Range("J120").Value = 35206
$ws.Range("L120").Value = 35206
$ws.Range("N120").Value = -44882
# Row 124 (Leve Item ID 34280)
$ws.Range("H124").Value = 34476.332
$ws.Range("J124").Value = 34476.332
$ws.Range("L124").Value = 34476.332
$ws.Range("N124").Value = -44296.332
# Row 128 (Leve Item ID 34563)
$ws.Range("H128").Value = 49715
$ws.Range("J128").Value = 49715
$ws.Range("L128").Value = 49715
$ws.Range("N128").Value = -59675
